$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 113, shifting existing rows 113:226 down to 114:227
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new data record
$ws.Cells.Item(113, 1).Value = 4
$ws.Cells.Item(113, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(113, 3).Value = "Los Lagos"
$ws.Cells.Item(113, 4).Value = 44638
$ws.Cells.Item(113, 5).Value = 10
$ws.Cells.Item(113, 6).Value = "Fruta"
$ws.Cells.Item(113, 7).Value = 100104
$ws.Cells.Item(113, 8).Value = "Frutos de pepita"
$ws.Cells.Item(113, 9).Value = 100104005
$ws.Cells.Item(113, 10).Value = "Pera"
$ws.Cells.Item(113, 11).Value = "Packham's Triumph"
$ws.Cells.Item(113, 12).Value = "Primera"
$ws.Cells.Item(113, 13).Value = 400
$ws.Cells.Item(113, 14).Value = 13000
$ws.Cells.Item(113, 15).Value = 14000
$ws.Cells.Item(113, 16).Value = 13500
$ws.Cells.Item(113, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(113, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(113, 19).Value = 900
$ws.Cells.Item(113, 20).Value = 15

# Ensure the date cell retains the date number format used by the rest of column D
$ws.Cells.Item(113, 4).NumberFormat = $ws.Cells.Item(114, 4).NumberFormat
